$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (including styles s=3 / s=5 and wrap text) from row 2 down to rows 3-5
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D5").PasteSpecial(-4122)

# Set the explicit custom row heights to match the new content
$ws.Rows.Item(3).RowHeight = 210
$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 195

# Set Date column (C) values for all four rows, in row order, so they occupy
# shared-string slots 6,7,8,9 respectively
$ws.Range("C2").Value = "09/09/2019"
$ws.Range("C3").Value = "09/10/2019"
$ws.Range("C4").Value = "09/11/2019"
$ws.Range("C5").Value = "09/12/2019"

# Set Remarks column (D) values from bottom row to top row, so they occupy
# shared-string slots 10,11,12,13 respectively (D5=10, D4=11, D3=12, D2=13)
$ws.Range("D5").Value = "1.	Worked on the test case update and test execution of #6642: Cart page redesign on Chrome browser. Created defects #8269 and #8270 as discussed. Will cancel defect #8269 as per discussion and create another defect from the observations discussed today. Please find updated test report attached.
2.	The dev preview is not working on dev environment. Shared the issue in a separate email.
3.	Attended daily stand-up to triage the defects in September release. The defects that are on dev for testing are marked as “On Dev” in smartsheet and highlighted in green. All the other defects that are in scope of September release are marked as “In Progress”  and highlighted in yellow. All the highlighted defects have been at the top under section “Defects” in September release.
4.	Attended the automation meeting with Noah. Discussed the new updates and queries related to Git and execution of the features ready to run."
$ws.Range("D4").Value = "1.	Added all open defects from August release to September release and consolidated the list under a single “Defects” section. Also, updated the status and developer’s name against all the defects.
2.	Worked on verification of the functional scenarios on the #6642: Cart page redesign and #8164: Contact us page add Google recaptcha. Also, updated test cases for ticket: #8164. Please find attached updated Test Report for September release."
$ws.Range("D3").Value = "1.	Verified the issue shared by you, with the space along the images appearing under two columns in section “Recommended For You” on mobile site on home page, on different versions of iOS and Android on different screen size and browsers through Browser-Stack. Shared analysis in another email.
2.	Worked on verification of the functional scenarios on the Cart page on desktop site as part of #6642. Also, updated test cases for ticket: #6642: Cart page redesign. Also, added the tickets deployed and tested last Friday in the September release test report. Please find the updated test report attached.
3.	Attended daily automation meeting with Noah and discussed queries and other updates.
4.	Updated the products/items against the regression test cases and shared updated sheet with Noah.
5.	Attended daily stand up meeting to discuss updates, queries, tasks and requirements."
$ws.Range("D2").Value = "1.	Worked on understanding of the Cucumber Framework elements (Features, Step Definitions, Hooks, Tags and Test Runner), implemented for automation as shown by Noah in our last call.  Also, analyzed the project created by Noah on Dev Server and worked on set- up of IJ IDE and other software’s required for Automation on my local. Set-up is still in progress.
2.	Attended the automation meeting with Noah and shared the challenges identified after the analyses of the BDD framework, implemented for personalization of the  products , so far. Let’s discuss the same in our daily stand up call tomorrow. Will continue with further analysis of the Automation project and share the feedback with you and Noah."

# Set Project (A) and Task (B) columns for the newly added rows (same values as row 2)
$ws.Range("A3").Value = "BBBY-PMall"
$ws.Range("B3").Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Range("A4").Value = "BBBY-PMall"
$ws.Range("B4").Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Range("A5").Value = "BBBY-PMall"
$ws.Range("B5").Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
